# Updates the "EC" worksheet's Valor Mora (column F) values for the
# periodo 2409 (row 16) and periodo 2404 (row 21) records, swapping
# the values that were previously entered on the wrong rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F16").Value = 52000
$ws.Range("F21").Value = 32933
